{"js": "// Replace the date line and the twenty-six two-digit multiplication\n// problems with their new values, as described by the diff.\nconst replacements = [\n  [\"2026-01-02 Friday\", \"2026-01-03 Saturday\"],\n  [\"41\u00d760=\", \"55\u00d736=\"],\n  [\"22\u00d751=\", \"66\u00d773=\"],\n  [\"89\u00d772=\", \"14\u00d773=\"],\n  [\"69\u00d778=\", \"22\u00d758=\"],\n  [\"15\u00d729=\", \"20\u00d717=\"],\n  [\"96\u00d733=\", \"95\u00d769=\"],\n  [\"41\u00d779=\", \"86\u00d716=\"],\n  [\"19\u00d766=\", \"23\u00d753=\"],\n  [\"53\u00d776=\", \"47\u00d772=\"],\n  [\"15\u00d716=\", \"18\u00d768=\"],\n  [\"85\u00d722=\", \"44\u00d776=\"],\n  [\"37\u00d747=\", \"88\u00d766=\"],\n  [\"65\u00d784=\", \"35\u00d746=\"],\n  [\"25\u00d726=\", \"20\u00d739=\"],\n  [\"36\u00d798=\", \"66\u00d718=\"],\n  [\"63\u00d786=\", \"35\u00d791=\"],\n  [\"21\u00d721=\", \"41\u00d766=\"],\n  [\"77\u00d739=\", \"53\u00d734=\"],\n  [\"82\u00d728=\", \"79\u00d778=\"],\n  [\"84\u00d790=\", \"19\u00d787=\"],\n  [\"58\u00d717=\", \"85\u00d792=\"],\n  [\"49\u00d771=\", \"23\u00d768=\"],\n  [\"49\u00d796=\", \"72\u00d772=\"],\n  [\"33\u00d755=\", \"76\u00d775=\"],\n  [\"50\u00d772=\", \"68\u00d739=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the twenty-six two-digit multiplication\n# problems with their new values, as described by the diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2026-01-02 Friday\", \"2026-01-03 Saturday\"),\n    @(\"41\u00d760=\", \"55\u00d736=\"),\n    @(\"22\u00d751=\", \"66\u00d773=\"),\n    @(\"89\u00d772=\", \"14\u00d773=\"),\n    @(\"69\u00d778=\", \"22\u00d758=\"),\n    @(\"15\u00d729=\", \"20\u00d717=\"),\n    @(\"96\u00d733=\", \"95\u00d769=\"),\n    @(\"41\u00d779=\", \"86\u00d716=\"),\n    @(\"19\u00d766=\", \"23\u00d753=\"),\n    @(\"53\u00d776=\", \"47\u00d772=\"),\n    @(\"15\u00d716=\", \"18\u00d768=\"),\n    @(\"85\u00d722=\", \"44\u00d776=\"),\n    @(\"37\u00d747=\", \"88\u00d766=\"),\n    @(\"65\u00d784=\", \"35\u00d746=\"),\n    @(\"25\u00d726=\", \"20\u00d739=\"),\n    @(\"36\u00d798=\", \"66\u00d718=\"),\n    @(\"63\u00d786=\", \"35\u00d791=\"),\n    @(\"21\u00d721=\", \"41\u00d766=\"),\n    @(\"77\u00d739=\", \"53\u00d734=\"),\n    @(\"82\u00d728=\", \"79\u00d778=\"),\n    @(\"84\u00d790=\", \"19\u00d787=\"),\n    @(\"58\u00d717=\", \"85\u00d792=\"),\n    @(\"49\u00d771=\", \"23\u00d768=\"),\n    @(\"49\u00d796=\", \"72\u00d772=\"),\n    @(\"33\u00d755=\", \"76\u00d775=\"),\n    @(\"50\u00d772=\", \"68\u00d739=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
